$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ODI Batting Extra")

# Insert 5 new rows right after the header row (row 1). This pushes the
# existing data rows (currently rows 2-21) down to rows 7-26, matching
# the new A1:F26 dimension.
$ws.Range("A2:F6").EntireRow.Insert()

# The inserted rows picked up the header row's bold/centered formatting;
# reset them back to the plain "Normal" style used by the other data rows.
$ws.Range("A2:F6").Style = "Normal"

# --- Row 2: MATCH_CODE=3460, BATTING_POSITION=9, NUM_4=0, NUM_6=0, MAN_OF_MATCH=NO ---
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "3460"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = 9

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "0"
$ws.Range("C2").Style = "Normal"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0"
$ws.Range("D2").Style = "Normal"

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "NO"
$ws.Range("F2").Style = "Normal"

# --- Row 3: MATCH_CODE=3659, BATTING_POSITION=10, MAN_OF_MATCH=NO ---
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "3659"
$ws.Range("A3").Style = "Normal"

$ws.Range("B3").Value = 10

$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "NO"
$ws.Range("F3").Style = "Normal"

# --- Row 4: MATCH_CODE=3663, BATTING_POSITION=11, MAN_OF_MATCH=NO ---
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "3663"
$ws.Range("A4").Style = "Normal"

$ws.Range("B4").Value = 11

$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "NO"
$ws.Range("F4").Style = "Normal"

# --- Row 5: MATCH_CODE=3669, BATTING_POSITION=10, MAN_OF_MATCH=NO ---
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "3669"
$ws.Range("A5").Style = "Normal"

$ws.Range("B5").Value = 10

$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "NO"
$ws.Range("F5").Style = "Normal"

# --- Row 6: MATCH_CODE=3677, BATTING_POSITION=(blank), MAN_OF_MATCH=NO ---
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "3677"
$ws.Range("A6").Style = "Normal"

$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "NO"
$ws.Range("F6").Style = "Normal"
